# Applies the 2023-05-03 daily crime-data update across all affected worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 10).Value = 2276
$ws.Cells.Item(3, 10).Value = 2369
$ws.Cells.Item(4, 8).Value = 1693
$ws.Cells.Item(4, 9).Value = 1757
$ws.Cells.Item(4, 10).Value = 540
$ws.Cells.Item(6, 10).Value = 2990
$ws.Cells.Item(7, 8).Value = 26006
$ws.Cells.Item(7, 9).Value = 26204
$ws.Cells.Item(7, 10).Value = 8344

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(6, 10).Value = 28
$ws.Cells.Item(7, 10).Value = 93

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 10).Value = 83
$ws.Cells.Item(3, 10).Value = 104
$ws.Cells.Item(6, 10).Value = 81
$ws.Cells.Item(7, 10).Value = 284

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(2, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 114

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 10).Value = 64
$ws.Cells.Item(3, 10).Value = 122
$ws.Cells.Item(7, 10).Value = 306

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(3, 10).Value = 16
$ws.Cells.Item(6, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 64

$ws = $wb.Worksheets.Item('South Deering')
$ws.Cells.Item(3, 10).Value = 19
$ws.Cells.Item(7, 10).Value = 79

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(3, 10).Value = 61
$ws.Cells.Item(7, 10).Value = 220

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 10).Value = 63
$ws.Cells.Item(7, 10).Value = 253
$ws.Cells.Item(8, 10).Value = 528
$ws.Cells.Item(15, 10).Value = 104
$ws.Cells.Item(18, 10).Value = 96
$ws.Cells.Item(19, 10).Value = 267
$ws.Cells.Item(20, 10).Value = 174
$ws.Cells.Item(24, 10).Value = 28
$ws.Cells.Item(29, 10).Value = 467
$ws.Cells.Item(31, 10).Value = 64
$ws.Cells.Item(33, 10).Value = 345
$ws.Cells.Item(35, 10).Value = 8
$ws.Cells.Item(36, 10).Value = 129
$ws.Cells.Item(37, 10).Value = 284
$ws.Cells.Item(40, 10).Value = 16
$ws.Cells.Item(44, 10).Value = 70
$ws.Cells.Item(47, 10).Value = 74
$ws.Cells.Item(48, 10).Value = 78
$ws.Cells.Item(51, 10).Value = 111
$ws.Cells.Item(52, 10).Value = 203
$ws.Cells.Item(54, 10).Value = 168
$ws.Cells.Item(63, 8).Value = 243
$ws.Cells.Item(63, 9).Value = 206
$ws.Cells.Item(63, 10).Value = 35
$ws.Cells.Item(65, 10).Value = 220
$ws.Cells.Item(66, 10).Value = 21
$ws.Cells.Item(67, 10).Value = 306
$ws.Cells.Item(72, 10).Value = 32
$ws.Cells.Item(76, 10).Value = 120
$ws.Cells.Item(78, 10).Value = 115
$ws.Cells.Item(79, 10).Value = 255
$ws.Cells.Item(80, 10).Value = 18
$ws.Cells.Item(83, 10).Value = 197
$ws.Cells.Item(84, 10).Value = 79
$ws.Cells.Item(85, 10).Value = 395
$ws.Cells.Item(86, 10).Value = 46
$ws.Cells.Item(90, 10).Value = 92
$ws.Cells.Item(93, 10).Value = 43
$ws.Cells.Item(94, 10).Value = 67
$ws.Cells.Item(96, 10).Value = 93
$ws.Cells.Item(98, 10).Value = 51
$ws.Cells.Item(99, 10).Value = 114
$ws.Cells.Item(101, 8).Value = 26006
$ws.Cells.Item(101, 9).Value = 26204
$ws.Cells.Item(101, 10).Value = 8344

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(3, 10).Value = 68
$ws.Cells.Item(7, 10).Value = 197

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 10).Value = 92
$ws.Cells.Item(4, 10).Value = 18
$ws.Cells.Item(7, 10).Value = 345

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(4, 10).Value = 12
$ws.Cells.Item(7, 10).Value = 168

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(3, 10).Value = 161
$ws.Cells.Item(6, 10).Value = 129
$ws.Cells.Item(7, 10).Value = 467

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(3, 10).Value = 74
$ws.Cells.Item(7, 10).Value = 267

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(7, 10).Value = 70

$ws = $wb.Worksheets.Item('Lake View')
$ws.Cells.Item(2, 10).Value = 15
$ws.Cells.Item(7, 10).Value = 78

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(2, 10).Value = 15
$ws.Cells.Item(4, 10).Value = 11
$ws.Cells.Item(7, 10).Value = 120

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 10).Value = 96
$ws.Cells.Item(3, 10).Value = 150
$ws.Cells.Item(4, 10).Value = 26
$ws.Cells.Item(7, 10).Value = 395

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 10).Value = 29
$ws.Cells.Item(7, 10).Value = 115

$ws = $wb.Worksheets.Item('Dunning')
$ws.Cells.Item(3, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 28

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 10).Value = 73
$ws.Cells.Item(3, 10).Value = 94
$ws.Cells.Item(6, 10).Value = 70
$ws.Cells.Item(7, 10).Value = 255

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 10).Value = 53
$ws.Cells.Item(4, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 174

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(6, 10).Value = 53
$ws.Cells.Item(7, 10).Value = 96

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(3, 10).Value = 31
$ws.Cells.Item(4, 10).Value = 3
$ws.Cells.Item(7, 10).Value = 129

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(2, 10).Value = 8
$ws.Cells.Item(7, 10).Value = 43

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(6, 10).Value = 85
$ws.Cells.Item(7, 10).Value = 203

$ws = $wb.Worksheets.Item('West Loop')
$ws.Cells.Item(2, 10).Value = 13
$ws.Cells.Item(6, 10).Value = 40
$ws.Cells.Item(7, 10).Value = 67

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(7, 10).Value = 74

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(3, 10).Value = 27
$ws.Cells.Item(7, 10).Value = 104

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(6, 10).Value = 28
$ws.Cells.Item(7, 10).Value = 51

$ws = $wb.Worksheets.Item('North Center')
$ws.Cells.Item(6, 10).Value = 13
$ws.Cells.Item(7, 10).Value = 21

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Cells.Item(6, 10).Value = 5
$ws.Cells.Item(7, 10).Value = 8

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(2, 10).Value = 22
$ws.Cells.Item(7, 10).Value = 63

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(2, 10).Value = 166
$ws.Cells.Item(6, 10).Value = 153
$ws.Cells.Item(7, 10).Value = 528

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(4, 10).Value = 21
$ws.Cells.Item(7, 10).Value = 46

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(7, 10).Value = 92

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(3, 10).Value = 32
$ws.Cells.Item(7, 10).Value = 111

$ws = $wb.Worksheets.Item('Old Town')
$ws.Cells.Item(2, 10).Value = 11
$ws.Cells.Item(6, 10).Value = 6
$ws.Cells.Item(7, 10).Value = 32

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Cells.Item(6, 10).Value = 12
$ws.Cells.Item(7, 10).Value = 18

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Cells.Item(3, 10).Value = 3
$ws.Cells.Item(7, 10).Value = 16

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 10).Value = 81
$ws.Cells.Item(6, 10).Value = 88
$ws.Cells.Item(7, 10).Value = 253

Write-Output "Applied 149 cell updates across 41 worksheets."
